$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.558.85'
$ws.Range('E2').Value = '  -2.09%  '
$ws.Range('D3').Value = '1.579.93'
$ws.Range('E3').Value = '  -3.04%  '
$ws.Range('E4').Value = '  +0.17%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '210.68'
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = '  -2.61%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.505'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -2.03%  '
$ws.Range('E7').Value = '  +0.18%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.248'
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = '  -2.04%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.0616'
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = '  -1.17%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '19.45'
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = '  -3.54%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0833'
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = '  -2.05%  '
$ws.Range('D12').Value = '1.802.50'
$ws.Range('E12').Value = '  -2.95%  '
$ws.Range('D13').Value = '1.571.93'
$ws.Range('E13').Value = '  -3.42%  '
$ws.Range('E14').Value = '  -1.50%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.526'
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = '  -2.54%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '63.73'
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = '  -1.47%  '
$ws.Range('D17').Value = '26.579.71'
$ws.Range('E17').Value = '  -1.90%  '
$ws.Range('D18').Value = '0.0₃0727'
$ws.Range('E18').Value = '  -0.66%  '
$ws.Range('E19').Value = '  +0.22%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '208.16'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  -2.63%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.66'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  -3.35%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.25'
$ws.Range('D22').Style = "Normal"
$ws.Range('E22').Value = '  -3.20%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '2.36'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  -5.59%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '8.87'
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = '  -1.85%  '
$ws.Range('E25').Value = '  -1.38%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '1.00'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '7.43'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  +2.22%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.113'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  -4.20%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '15.25'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  -1.77%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '0.0499'
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = '  -0.81%  '
$ws.Range('E31').Value = '  -2.21%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.24'
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = '  -4.06%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.655'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +22.35%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '2.93'
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = '  -2.58%  '
$ws.Range('D35').Value = '1.305.60'
$ws.Range('E35').Value = '  -0.79%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.51'
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = '  -3.37%  '
$ws.Range('E37').Value = '  -1.02%  '
$ws.Range('E38').Value = '  -1.16%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.818'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -2.86%  '
$ws.Range('E40').Value = '  +0.20%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.786'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  -2.07%  '
$ws.Range('B42').Value = 'FraxShare'
$ws.Range('C42').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.28'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +0.71%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.16'
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = '  -4.37%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '62.64'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  -2.13%  '
$ws.Range('D45').Value = '1.716.65'
$ws.Range('E45').Value = '  -2.69%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '88.81'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  -2.02%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.61'
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = '  +0.62%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '0.829'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +5.09%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '0.0506'
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = '  -1.72%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.0980'
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = '  +3.73%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₇0964'
$ws.Range('E51').Value = '  -8.90%  '
